# Update "想去人数" (F column) values on the sheets that hold the full
# data table: "展览" (sheet 1) and "全部类型" (sheet 4). Sheets "演出" and
# "本地生活" only contain a header row, so nothing changes there.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1752
    $ws.Range("F3").Value = 8026
    $ws.Range("F4").Value = 186
    $ws.Range("F5").Value = 281
}
